$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Candidate ID (numeric)
$ws.Range("B2").Value = 23081034

# Update User Name
$ws.Range("C2").Value = "ubzzrnn59"

# Update Exam Password
$ws.Range("D2").Value = "rj82JD&$"

# Update First Name
$ws.Range("F2").Value = "gknOSSlg"

# Update Last Name
$ws.Range("G2").Value = "LCJs"

# Update Client Id
$ws.Range("A2").Value = "gATGO674"
